# Add a new "Prueba" material row (row 17) to the materials table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B17").Value = 15
$ws.Range("C17").Value = "Prueba"
$ws.Range("D17").Value = 2888
$ws.Range("E17").Value = 200000000000
$ws.Range("F17").Value = 0.0001
$ws.Range("G17").Value = 0.12

# Match the author's final selection/scroll state after entering the row.
$ws.Range("A17:F17").Select()
$excel.ActiveWindow.ScrollRow = 7
